# Ran code for averaged intensities on spiral schemes
#
# The averaged-intensity table gains four new sampling schemes
# (Gaussian-Quadrature, relocated, plus three new Spiral-* schemes).
# They are inserted right after "Ring Perpendicular to TD" (pushing the
# NoRotation/Rotation/HexGrid rows down), and three brand new rows are
# appended at the bottom of the table to make room for the three
# HexGrid schemes that no longer fit within the original 14-row block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ordered list of sampling-scheme names for rows 10..19 (column B);
# rows 3..9 (ND Single .. Ring Perpendicular to TD) are unchanged.
# Column A holds the 1-based scheme index.
$schemes = @(
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)

$cols = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16)

for ($i = 0; $i -lt $schemes.Length; $i++) {
    $row = $i + 10

    if ($row -le 16) {
        # Existing row: only the scheme index/name change, formatting stays.
        $ws.Cells.Item($row, 1).Value = $i + 8
        $ws.Cells.Item($row, 2).Value = $schemes[$i]
    } else {
        # Brand new row: clone formatting from the row above, then fill in.
        $ws.Cells.Item($row, 1).Value = $i + 8
        $ws.Cells.Item($row - 1, 1).Copy()
        $ws.Cells.Item($row, 1).PasteSpecial(-4122)
        $ws.Cells.Item($row, 2).Value = $schemes[$i]
        foreach ($c in $cols) {
            $ws.Cells.Item($row, $c).Value = 1
        }
    }
}
